# Commit: "se vuelve a 6 de entrada" (back to 6 on entry)
# - On sheet "grilla de pruebas": B9 goes back to 6 (previously 9), and loses the
#   yellow highlight / text-format styling it had, taking on the same plain style
#   used by its neighboring input cells (e.g. B8).
# - The accompanying note in C9 is reworded.
# - Selection on that sheet moves to C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grilla de pruebas")

$B9 = $ws.Range("B9")
$B8 = $ws.Range("B8")

# Set the numeric value first so the paste-formats step (which brings in B8's
# "@" text number-format) does not coerce it into a text string.
$B9.Value = 6

# Copy just the formatting from B8 (no fill / thin border / text number format)
# onto B9 so it matches its neighbours again instead of the old yellow highlight.
$B8.Copy()
$B9.PasteSpecial(-4122)

# Update the note text next to it.
$ws.Range("C9").Value = "Automáico 6. manual 19"

# Move the active selection to C10 on this sheet.
$ws.Activate()
$ws.Range("C10").Select()
